# Auto-generated: apply cryptos.xlsx data-refresh diff (Mon May  8 09:23:06 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.896.83'
$ws.Range('E2').Value = '  -3.46%  '
$ws.Range('D3').Value = '1.856.58'
$ws.Range('E3').Value = '  -2.73%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.85'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4344'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3678'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.73%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07484'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9376'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.32'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.33%  '
$ws.Range('D12').Value = '1.832.12'
$ws.Range('E12').Value = '  -3.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.691'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.427'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.21%  '
$ws.Range('E15').Value = '  -2.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.003'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '81.36'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009034'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.003'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('E20').Value = '  -4.79%  '
$ws.Range('D21').Value = '27.885.65'
$ws.Range('E21').Value = '  -3.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.098'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.93'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('D24').Value = '2.099.05'
$ws.Range('E24').Value = '  -2.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.004'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.40%  '
$ws.Range('E26').Value = '  -2.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.32'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.395'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.26'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.737'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08971'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.38%  '
$ws.Range('E32').Value = '  -7.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.819'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.999'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.169'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.003'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.113'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05414'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01962'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.920'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5237'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.999'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1679'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.26%  '
$ws.Range('E44').Value = '  -6.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.06703'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.85%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4871'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.61'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '106.77'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.920'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -8.10%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.675'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.73%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.09%  '
